# ------------------------------------------------------------------------
# Refresh the coinranking.com snapshot table on Sheet1 (rows 2-51) to the
# Wed Jan 11 13:26:03 UTC 2023 GitHub Actions run:
#   - column G ("Hora") moves from 12 -> 13 for every data row
#   - columns D ("Price") / E ("Volume(1h)") are refreshed with the latest
#     quote for rows whose values actually moved
#   - rows 10-24 absorb one more coin from the ranking feed, so each row
#     B (Coin) / C (Link) pair shifts up to the row above it, with its own
#     refreshed D/E figures
#
# All of these source cells are stored as literal text (e.g. "1.42%", not
# a 1.42% formatted number), so every write below is bracketed with a
# "@" (Text) number format and a ClearFormats() to land back on the sheet`s
# original (unstyled) cell format while keeping the value a text string.
# ------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:G2").NumberFormat = "@"
$ws.Range("D2").Value = "278.72"
$ws.Range("E2").Value = "1.42%"
$ws.Range("G2").Value = "13"
$ws.Range("D2:G2").ClearFormats()

# Row 3
$ws.Range("D3:G3").NumberFormat = "@"
$ws.Range("D3").Value = "27.22"
$ws.Range("E3").Value = "2.22%"
$ws.Range("G3").Value = "13"
$ws.Range("D3:G3").ClearFormats()

# Row 4
$ws.Range("D4:G4").NumberFormat = "@"
$ws.Range("D4").Value = "4.855"
$ws.Range("E4").Value = "2.00%"
$ws.Range("G4").Value = "13"
$ws.Range("D4:G4").ClearFormats()

# Row 5
$ws.Range("D5:G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06409"
$ws.Range("E5").Value = "1.88%"
$ws.Range("G5").Value = "13"
$ws.Range("D5:G5").ClearFormats()

# Row 6
$ws.Range("D6:G6").NumberFormat = "@"
$ws.Range("D6").Value = "6.996"
$ws.Range("E6").Value = "0.98%"
$ws.Range("G6").Value = "13"
$ws.Range("D6:G6").ClearFormats()

# Row 7
$ws.Range("D7:G7").NumberFormat = "@"
$ws.Range("D7").Value = "1.222"
$ws.Range("E7").Value = "-8.49%"
$ws.Range("G7").Value = "13"
$ws.Range("D7:G7").ClearFormats()

# Row 8
$ws.Range("D8:G8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8818"
$ws.Range("E8").Value = "1.38%"
$ws.Range("G8").Value = "13"
$ws.Range("D8:G8").ClearFormats()

# Row 9
$ws.Range("D9:G9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1528"
$ws.Range("E9").Value = "-2.15%"
$ws.Range("G9").Value = "13"
$ws.Range("D9:G9").ClearFormats()

# Row 10
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10:G10").NumberFormat = "@"
$ws.Range("D10").Value = "0.05184"
$ws.Range("E10").Value = "3.11%"
$ws.Range("G10").Value = "13"
$ws.Range("D10:G10").ClearFormats()

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11:G11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07520"
$ws.Range("E11").Value = "0.56%"
$ws.Range("G11").Value = "13"
$ws.Range("D11:G11").ClearFormats()

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12:G12").NumberFormat = "@"
$ws.Range("D12").Value = "0.02873"
$ws.Range("E12").Value = "-2.54%"
$ws.Range("G12").Value = "13"
$ws.Range("D12:G12").ClearFormats()

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13:G13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08962"
$ws.Range("E13").Value = "-1.00%"
$ws.Range("G13").Value = "13"
$ws.Range("D13:G13").ClearFormats()

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14:G14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001563"
$ws.Range("E14").Value = "-0.07%"
$ws.Range("G14").Value = "13"
$ws.Range("D14:G14").ClearFormats()

# Row 15
$ws.Range("B15").Value = "One"
$ws.Range("C15").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D15:G15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0006363"
$ws.Range("E15").Value = "0.80%"
$ws.Range("G15").Value = "13"
$ws.Range("D15:G15").ClearFormats()

# Row 16
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16:G16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006092"
$ws.Range("E16").Value = "2.04%"
$ws.Range("G16").Value = "13"
$ws.Range("D16:G16").ClearFormats()

# Row 17
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17:G17").NumberFormat = "@"
$ws.Range("D17").Value = "3.480"
$ws.Range("E17").Value = "0.78%"
$ws.Range("G17").Value = "13"
$ws.Range("D17:G17").ClearFormats()

# Row 18
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18:G18").NumberFormat = "@"
$ws.Range("D18").Value = "3.302"
$ws.Range("E18").Value = "-0.19%"
$ws.Range("G18").Value = "13"
$ws.Range("D18:G18").ClearFormats()

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19:G19").NumberFormat = "@"
$ws.Range("D19").Value = "2.247"
$ws.Range("E19").Value = "-1.60%"
$ws.Range("G19").Value = "13"
$ws.Range("D19:G19").ClearFormats()

# Row 20
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20:G20").NumberFormat = "@"
$ws.Range("D20").Value = "0.3085"
$ws.Range("E20").Value = "-2.03%"
$ws.Range("G20").Value = "13"
$ws.Range("D20:G20").ClearFormats()

# Row 21
$ws.Range("E21:G21").NumberFormat = "@"
$ws.Range("E21").Value = "2.36%"
$ws.Range("G21").Value = "13"
$ws.Range("E21:G21").ClearFormats()

# Row 22
$ws.Range("D22:G22").NumberFormat = "@"
$ws.Range("D22").Value = "3.903"
$ws.Range("E22").Value = "-0.35%"
$ws.Range("G22").Value = "13"
$ws.Range("D22:G22").ClearFormats()

# Row 23
$ws.Range("B23").Value = "CoinExToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D23:G23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04401"
$ws.Range("E23").Value = "0.69%"
$ws.Range("G23").Value = "13"
$ws.Range("D23:G23").ClearFormats()

# Row 24
$ws.Range("B24").Value = "ZBToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D24:G24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1504"
$ws.Range("E24").Value = "8.95%"
$ws.Range("G24").Value = "13"
$ws.Range("D24:G24").ClearFormats()

# Row 25
$ws.Range("E25:G25").NumberFormat = "@"
$ws.Range("E25").Value = "0.54%"
$ws.Range("G25").Value = "13"
$ws.Range("E25:G25").ClearFormats()

# Row 26
$ws.Range("D26:G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.003900"
$ws.Range("E26").Value = "-7.36%"
$ws.Range("G26").Value = "13"
$ws.Range("D26:G26").ClearFormats()

# Row 27
$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "13"
$ws.Range("G27").ClearFormats()

# Row 28
$ws.Range("E28:G28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.70%"
$ws.Range("G28").Value = "13"
$ws.Range("E28:G28").ClearFormats()

# Row 29
$ws.Range("D29:G29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0001642"
$ws.Range("E29").Value = "1.57%"
$ws.Range("G29").Value = "13"
$ws.Range("D29:G29").ClearFormats()

# Row 30
$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "13"
$ws.Range("G30").ClearFormats()

# Row 31
$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "13"
$ws.Range("G31").ClearFormats()

# Row 32
$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "13"
$ws.Range("G32").ClearFormats()

# Row 33
$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "13"
$ws.Range("G33").ClearFormats()

# Row 34
$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "13"
$ws.Range("G34").ClearFormats()

# Row 35
$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "13"
$ws.Range("G35").ClearFormats()

# Row 36
$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "13"
$ws.Range("G36").ClearFormats()

# Row 37
$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "13"
$ws.Range("G37").ClearFormats()

# Row 38
$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "13"
$ws.Range("G38").ClearFormats()

# Row 39
$ws.Range("G39").NumberFormat = "@"
$ws.Range("G39").Value = "13"
$ws.Range("G39").ClearFormats()

# Row 40
$ws.Range("D40:G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04091"
$ws.Range("E40").Value = "0.58%"
$ws.Range("G40").Value = "13"
$ws.Range("D40:G40").ClearFormats()

# Row 41
$ws.Range("D41:G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006826"
$ws.Range("E41").Value = "-3.16%"
$ws.Range("G41").Value = "13"
$ws.Range("D41:G41").ClearFormats()

# Row 42
$ws.Range("E42:G42").NumberFormat = "@"
$ws.Range("E42").Value = "0.48%"
$ws.Range("G42").Value = "13"
$ws.Range("E42:G42").ClearFormats()

# Row 43
$ws.Range("E43:G43").NumberFormat = "@"
$ws.Range("E43").Value = "-5.54%"
$ws.Range("G43").Value = "13"
$ws.Range("E43:G43").ClearFormats()

# Row 44
$ws.Range("D44:G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01125"
$ws.Range("E44").Value = "0.16%"
$ws.Range("G44").Value = "13"
$ws.Range("D44:G44").ClearFormats()

# Row 45
$ws.Range("D45:G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005345"
$ws.Range("E45").Value = "2.68%"
$ws.Range("G45").Value = "13"
$ws.Range("D45:G45").ClearFormats()

# Row 46
$ws.Range("E46:G46").NumberFormat = "@"
$ws.Range("E46").Value = "9.53%"
$ws.Range("G46").Value = "13"
$ws.Range("E46:G46").ClearFormats()

# Row 47
$ws.Range("D47:G47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01850"
$ws.Range("E47").Value = "-19.66%"
$ws.Range("G47").Value = "13"
$ws.Range("D47:G47").ClearFormats()

# Row 48
$ws.Range("G48").NumberFormat = "@"
$ws.Range("G48").Value = "13"
$ws.Range("G48").ClearFormats()

# Row 49
$ws.Range("G49").NumberFormat = "@"
$ws.Range("G49").Value = "13"
$ws.Range("G49").ClearFormats()

# Row 50
$ws.Range("G50").NumberFormat = "@"
$ws.Range("G50").Value = "13"
$ws.Range("G50").ClearFormats()

# Row 51
$ws.Range("G51").NumberFormat = "@"
$ws.Range("G51").Value = "13"
$ws.Range("G51").ClearFormats()
